$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-15 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-16 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("37-9=28", $true, $false, $false, $false, $false, $true, 1, $false, "45+33=78", 2) | Out-Null
$d.Content.Find.Execute("3+80=83", $true, $false, $false, $false, $false, $true, 1, $false, "38+50=88", 2) | Out-Null
$d.Content.Find.Execute("98+1=99", $true, $false, $false, $false, $false, $true, 1, $false, "33-10=23", 2) | Out-Null
$d.Content.Find.Execute("30+65=95", $true, $false, $false, $false, $false, $true, 1, $false, "47-21=26", 2) | Out-Null
$d.Content.Find.Execute("18+69=87", $true, $false, $false, $false, $false, $true, 1, $false, "50+32=82", 2) | Out-Null
$d.Content.Find.Execute("50+49=99", $true, $false, $false, $false, $false, $true, 1, $false, "26+42=68", 2) | Out-Null
$d.Content.Find.Execute("24-17=7", $true, $false, $false, $false, $false, $true, 1, $false, "66+20=86", 2) | Out-Null
$d.Content.Find.Execute("68+2=70", $true, $false, $false, $false, $false, $true, 1, $false, "21+23=44", 2) | Out-Null
$d.Content.Find.Execute("13+33=46", $true, $false, $false, $false, $false, $true, 1, $false, "53-46=7", 2) | Out-Null
$d.Content.Find.Execute("54+23=77", $true, $false, $false, $false, $false, $true, 1, $false, "16-8=8", 2) | Out-Null
$d.Content.Find.Execute("20+52=72", $true, $false, $false, $false, $false, $true, 1, $false, "20-17=3", 2) | Out-Null
$d.Content.Find.Execute("36+16=52", $true, $false, $false, $false, $false, $true, 1, $false, "84-15=69", 2) | Out-Null
$d.Content.Find.Execute("93-75=18", $true, $false, $false, $false, $false, $true, 1, $false, "46-42=4", 2) | Out-Null
$d.Content.Find.Execute("98-1=97", $true, $false, $false, $false, $false, $true, 1, $false, "54-22=32", 2) | Out-Null
$d.Content.Find.Execute("38-15=23", $true, $false, $false, $false, $false, $true, 1, $false, "55+18=73", 2) | Out-Null
$d.Content.Find.Execute("97-85=12", $true, $false, $false, $false, $false, $true, 1, $false, "69+9=78", 2) | Out-Null
$d.Content.Find.Execute("2+60=62", $true, $false, $false, $false, $false, $true, 1, $false, "7+78=85", 2) | Out-Null
$d.Content.Find.Execute("55+20=75", $true, $false, $false, $false, $false, $true, 1, $false, "55+21=76", 2) | Out-Null
$d.Content.Find.Execute("86-62=24", $true, $false, $false, $false, $false, $true, 1, $false, "57-18=39", 2) | Out-Null
$d.Content.Find.Execute("73+24=97", $true, $false, $false, $false, $false, $true, 1, $false, "11-11=0", 2) | Out-Null
$d.Content.Find.Execute("12+31=43", $true, $false, $false, $false, $false, $true, 1, $false, "14+61=75", 2) | Out-Null
$d.Content.Find.Execute("71+7=78", $true, $false, $false, $false, $false, $true, 1, $false, "20+74=94", 2) | Out-Null
$d.Content.Find.Execute("77-46=31", $true, $false, $false, $false, $false, $true, 1, $false, "12+15=27", 2) | Out-Null
$d.Content.Find.Execute("17+62=79", $true, $false, $false, $false, $false, $true, 1, $false, "57+31=88", 2) | Out-Null
$d.Content.Find.Execute("37-29=8", $true, $false, $false, $false, $false, $true, 1, $false, "83-0=83", 2) | Out-Null
$d.Content.Find.Execute("58-18=40", $true, $false, $false, $false, $false, $true, 1, $false, "65-29=36", 2) | Out-Null
$d.Content.Find.Execute("95-27=68", $true, $false, $false, $false, $false, $true, 1, $false, "31+17=48", 2) | Out-Null
$d.Content.Find.Execute("2+59=61", $true, $false, $false, $false, $false, $true, 1, $false, "22+1=23", 2) | Out-Null
$d.Content.Find.Execute("27+58=85", $true, $false, $false, $false, $false, $true, 1, $false, "47+6=53", 2) | Out-Null
$d.Content.Find.Execute("24+31=55", $true, $false, $false, $false, $false, $true, 1, $false, "52-40=12", 2) | Out-Null
$d.Content.Find.Execute("75-48=27", $true, $false, $false, $false, $false, $true, 1, $false, "13+6=19", 2) | Out-Null
$d.Content.Find.Execute("93+5=98", $true, $false, $false, $false, $false, $true, 1, $false, "56+7=63", 2) | Out-Null
$d.Content.Find.Execute("65-47=18", $true, $false, $false, $false, $false, $true, 1, $false, "80-35=45", 2) | Out-Null
$d.Content.Find.Execute("89-58=31", $true, $false, $false, $false, $false, $true, 1, $false, "18+15=33", 2) | Out-Null
$d.Content.Find.Execute("74-68=6", $true, $false, $false, $false, $false, $true, 1, $false, "60+1=61", 2) | Out-Null
$d.Content.Find.Execute("28+28=56", $true, $false, $false, $false, $false, $true, 1, $false, "34+49=83", 2) | Out-Null
$d.Content.Find.Execute("40+39=79", $true, $false, $false, $false, $false, $true, 1, $false, "42+41=83", 2) | Out-Null
$d.Content.Find.Execute("47-29=18", $true, $false, $false, $false, $false, $true, 1, $false, "93-4=89", 2) | Out-Null
$d.Content.Find.Execute("8-5=3", $true, $false, $false, $false, $false, $true, 1, $false, "20+58=78", 2) | Out-Null
$d.Content.Find.Execute("58-39=19", $true, $false, $false, $false, $false, $true, 1, $false, "72+9=81", 2) | Out-Null
$d.Content.Find.Execute("57+12=69", $true, $false, $false, $false, $false, $true, 1, $false, "3+47=50", 2) | Out-Null
$d.Content.Find.Execute("41+10=51", $true, $false, $false, $false, $false, $true, 1, $false, "89-5=84", 2) | Out-Null
$d.Content.Find.Execute("84+2=86", $true, $false, $false, $false, $false, $true, 1, $false, "31+25=56", 2) | Out-Null
$d.Content.Find.Execute("29+16=45", $true, $false, $false, $false, $false, $true, 1, $false, "47-17=30", 2) | Out-Null
$d.Content.Find.Execute("22+7=29", $true, $false, $false, $false, $false, $true, 1, $false, "78-23=55", 2) | Out-Null
$d.Content.Find.Execute("35+63=98", $true, $false, $false, $false, $false, $true, 1, $false, "64+29=93", 2) | Out-Null
$d.Content.Find.Execute("7+11=18", $true, $false, $false, $false, $false, $true, 1, $false, "63+13=76", 2) | Out-Null
$d.Content.Find.Execute("43+29=72", $true, $false, $false, $false, $false, $true, 1, $false, "16+35=51", 2) | Out-Null
$d.Content.Find.Execute("98-21=77", $true, $false, $false, $false, $false, $true, 1, $false, "60-53=7", 2) | Out-Null
$d.Content.Find.Execute("90-5=85", $true, $false, $false, $false, $false, $true, 1, $false, "79-75=4", 2) | Out-Null
$d.Content.Find.Execute("27-2=25", $true, $false, $false, $false, $false, $true, 1, $false, "53+43=96", 2) | Out-Null
$d.Content.Find.Execute("36-8=28", $true, $false, $false, $false, $false, $true, 1, $false, "45+37=82", 2) | Out-Null
$d.Content.Find.Execute("42+25=67", $true, $false, $false, $false, $false, $true, 1, $false, "47-37=10", 2) | Out-Null
$d.Content.Find.Execute("98-3=95", $true, $false, $false, $false, $false, $true, 1, $false, "91-78=13", 2) | Out-Null
$d.Content.Find.Execute("36+27=63", $true, $false, $false, $false, $false, $true, 1, $false, "61+15=76", 2) | Out-Null
$d.Content.Find.Execute("19+53=72", $true, $false, $false, $false, $false, $true, 1, $false, "12+23=35", 2) | Out-Null
$d.Content.Find.Execute("10+78=88", $true, $false, $false, $false, $false, $true, 1, $false, "40+19=59", 2) | Out-Null
$d.Content.Find.Execute("44+11=55", $true, $false, $false, $false, $false, $true, 1, $false, "79-53=26", 2) | Out-Null
$d.Content.Find.Execute("83-59=24", $true, $false, $false, $false, $false, $true, 1, $false, "39-19=20", 2) | Out-Null
$d.Content.Find.Execute("64-41=23", $true, $false, $false, $false, $false, $true, 1, $false, "13-0=13", 2) | Out-Null
$d.Content.Find.Execute("56-1=55", $true, $false, $false, $false, $false, $true, 1, $false, "96-44=52", 2) | Out-Null
$d.Content.Find.Execute("51-29=22", $true, $false, $false, $false, $false, $true, 1, $false, "55+24=79", 2) | Out-Null
$d.Content.Find.Execute("55-46=9", $true, $false, $false, $false, $false, $true, 1, $false, "79-74=5", 2) | Out-Null
$d.Content.Find.Execute("96-45=51", $true, $false, $false, $false, $false, $true, 1, $false, "8+50=58", 2) | Out-Null
$d.Content.Find.Execute("58-20=38", $true, $false, $false, $false, $false, $true, 1, $false, "15+77=92", 2) | Out-Null
$d.Content.Find.Execute("17+25=42", $true, $false, $false, $false, $false, $true, 1, $false, "72-22=50", 2) | Out-Null
$d.Content.Find.Execute("93+0=93", $true, $false, $false, $false, $false, $true, 1, $false, "21-8=13", 2) | Out-Null
$d.Content.Find.Execute("25+29=54", $true, $false, $false, $false, $false, $true, 1, $false, "59-19=40", 2) | Out-Null
$d.Content.Find.Execute("68-43=25", $true, $false, $false, $false, $false, $true, 1, $false, "99-73=26", 2) | Out-Null
$d.Content.Find.Execute("19+79=98", $true, $false, $false, $false, $false, $true, 1, $false, "29+21=50", 2) | Out-Null
$d.Content.Find.Execute("27+27=54", $true, $false, $false, $false, $false, $true, 1, $false, "38+45=83", 2) | Out-Null
$d.Content.Find.Execute("66+29=95", $true, $false, $false, $false, $false, $true, 1, $false, "16+13=29", 2) | Out-Null
$d.Content.Find.Execute("39-37=2", $true, $false, $false, $false, $false, $true, 1, $false, "50-37=13", 2) | Out-Null
$d.Content.Find.Execute("55+2=57", $true, $false, $false, $false, $false, $true, 1, $false, "64-56=8", 2) | Out-Null
$d.Content.Find.Execute("58-42=16", $true, $false, $false, $false, $false, $true, 1, $false, "38+17=55", 2) | Out-Null
$d.Content.Find.Execute("71-36=35", $true, $false, $false, $false, $false, $true, 1, $false, "52-45=7", 2) | Out-Null
$d.Content.Find.Execute("69-32=37", $true, $false, $false, $false, $false, $true, 1, $false, "57+41=98", 2) | Out-Null
$d.Content.Find.Execute("52-46=6", $true, $false, $false, $false, $false, $true, 1, $false, "5+82=87", 2) | Out-Null
$d.Content.Find.Execute("78-15=63", $true, $false, $false, $false, $false, $true, 1, $false, "96-46=50", 2) | Out-Null
$d.Content.Find.Execute("93-77=16", $true, $false, $false, $false, $false, $true, 1, $false, "83-52=31", 2) | Out-Null
$d.Content.Find.Execute("60+24=84", $true, $false, $false, $false, $false, $true, 1, $false, "71-37=34", 2) | Out-Null
$d.Content.Find.Execute("11+87=98", $true, $false, $false, $false, $false, $true, 1, $false, "2+97=99", 2) | Out-Null
$d.Content.Find.Execute("69+3=72", $true, $false, $false, $false, $false, $true, 1, $false, "58+14=72", 2) | Out-Null
$d.Content.Find.Execute("67-56=11", $true, $false, $false, $false, $false, $true, 1, $false, "93-43=50", 2) | Out-Null
$d.Content.Find.Execute("96-17=79", $true, $false, $false, $false, $false, $true, 1, $false, "79-49=30", 2) | Out-Null
$d.Content.Find.Execute("21-4=17", $true, $false, $false, $false, $false, $true, 1, $false, "62-34=28", 2) | Out-Null
$d.Content.Find.Execute("19+38=57", $true, $false, $false, $false, $false, $true, 1, $false, "92-63=29", 2) | Out-Null
$d.Content.Find.Execute("58-44=14", $true, $false, $false, $false, $false, $true, 1, $false, "72-9=63", 2) | Out-Null
$d.Content.Find.Execute("56-41=15", $true, $false, $false, $false, $false, $true, 1, $false, "92-1=91", 2) | Out-Null
$d.Content.Find.Execute("16-7=9", $true, $false, $false, $false, $false, $true, 1, $false, "43-18=25", 2) | Out-Null
$d.Content.Find.Execute("67+6=73", $true, $false, $false, $false, $false, $true, 1, $false, "42+5=47", 2) | Out-Null
$d.Content.Find.Execute("24+12=36", $true, $false, $false, $false, $false, $true, 1, $false, "26+18=44", 2) | Out-Null
$d.Content.Find.Execute("17+41=58", $true, $false, $false, $false, $false, $true, 1, $false, "74-16=58", 2) | Out-Null
$d.Content.Find.Execute("80+5=85", $true, $false, $false, $false, $false, $true, 1, $false, "69-45=24", 2) | Out-Null
$d.Content.Find.Execute("65+32=97", $true, $false, $false, $false, $false, $true, 1, $false, "12+81=93", 2) | Out-Null
$d.Content.Find.Execute("44+43=87", $true, $false, $false, $false, $false, $true, 1, $false, "12+30=42", 2) | Out-Null
$d.Content.Find.Execute("9+69=78", $true, $false, $false, $false, $false, $true, 1, $false, "16+50=66", 2) | Out-Null
$d.Content.Find.Execute("73+21=94", $true, $false, $false, $false, $false, $true, 1, $false, "71-23=48", 2) | Out-Null
$d.Content.Find.Execute("82+7=89", $true, $false, $false, $false, $false, $true, 1, $false, "35+28=63", 2) | Out-Null
$d.Content.Find.Execute("29+66=95", $true, $false, $false, $false, $false, $true, 1, $false, "80-23=57", 2) | Out-Null
